# Apply cryptos price/volume updates published Wed May 15 13:48:30 UTC 2024 (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.378.98"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "2.970.27"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'581.66"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'152.55"
$ws.Range("E6").Value = "  +4.94%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "2.968.84"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'0.446"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'34.35"
$ws.Range("E14").Value = "  +5.01%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "3.463.92"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "64.333.69"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("E18").Value = "  +3.83%  "
$ws.Range("D19").Value = "2.973.62"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "'453.66"
$ws.Range("E20").Value = "  +4.33%  "
$ws.Range("D21").Value = "'13.58"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").Value = "'0.673"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("D24").Value = "'80.50"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "'10.91"
$ws.Range("E25").Value = "  +7.20%  "
$ws.Range("D26").Value = "'12.28"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Value = "'2.18"
$ws.Range("E27").Value = "  +6.38%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'7.62"
$ws.Range("E29").Value = "  +7.79%  "
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").Value = "'26.57"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.974"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  +8.62%  "
$ws.Range("D38").Value = "'5.61"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'49.04"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'43.88"
$ws.Range("E41").Value = "  +14.17%  "
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").Value = "'0.291"
$ws.Range("E43").Value = "  +7.91%  "
$ws.Range("D44").Value = "'8.31"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'378.91"
$ws.Range("E45").Value = "  +10.24%  "
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("D47").Value = "2.755.40"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").Value = "'134.66"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("E51").Value = "  +6.85%  "
